$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R values, mirroring the existing date/number series in column Q.
$ws.Range("R1").Value = 44066
$ws.Range("R2").Value = 1.7999999999999998
$ws.Range("R3").Value = 0.6
$ws.Range("R4").Value = 3.125
$ws.Range("R5").Value = 4.0750000000000002
$ws.Range("R6").Value = 6.1499999999999995
$ws.Range("R7").Value = 1.9000000000000001
$ws.Range("R8").Value = 6.4250000000000007
$ws.Range("R9").Value = 2.9250000000000003
$ws.Range("R10").Value = 0.6
$ws.Range("R11").Value = 1.55
$ws.Range("R12").Value = 2.8499999999999996

# Match the formatting of the adjoining column Q for the new column.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("Q2:Q12").Copy()
$ws.Range("R2:R12").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("A13").Select()
